$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.066.37"
$ws.Range("E2").Value = "  -4.37%  "
$ws.Range("D3").Value = "3.528.92"
$ws.Range("E3").Value = "  -5.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'571.74"
$ws.Range("E5").Value = "  -6.83%  "
$ws.Range("D6").Value = "'188.08"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("D7").Value = "3.526.48"
$ws.Range("E7").Value = "  -5.14%  "
$ws.Range("D8").Value = "'0.608"
$ws.Range("E8").Value = "  -4.68%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "'0.661"
$ws.Range("E10").Value = "  -8.08%  "
$ws.Range("D11").Value = "'0.143"
$ws.Range("E11").Value = "  -10.63%  "
$ws.Range("D12").Value = "'52.38"
$ws.Range("E12").Value = "  -9.76%  "
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  -12.50%  "
$ws.Range("D14").Value = "'9.68"
$ws.Range("E14").Value = "  -8.94%  "
$ws.Range("D15").Value = "4.093.50"
$ws.Range("E15").Value = "  -5.08%  "
$ws.Range("D16").Value = "3.529.70"
$ws.Range("E16").Value = "  -5.09%  "
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "'18.13"
$ws.Range("E18").Value = "  -6.30%  "
$ws.Range("D19").Value = "65.844.89"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "'12.00"
$ws.Range("E20").Value = "  -7.43%  "
$ws.Range("E21").Value = "  -8.25%  "
$ws.Range("D22").Value = "'389.05"
$ws.Range("E22").Value = "  -5.62%  "
$ws.Range("D23").Value = "'4.24"
$ws.Range("E23").Value = "  -7.71%  "
$ws.Range("D24").Value = "'84.67"
$ws.Range("E24").Value = "  -5.51%  "
$ws.Range("D25").Value = "'10.93"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'2.85"
$ws.Range("E26").Value = "  -6.62%  "
$ws.Range("D27").Value = "'12.22"
$ws.Range("E27").Value = "  -5.40%  "
$ws.Range("D28").Value = "'6.03"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  -7.99%  "
$ws.Range("D30").Value = "'8.79"
$ws.Range("E30").Value = "  -9.23%  "
$ws.Range("D31").Value = "'30.69"
$ws.Range("E31").Value = "  -7.50%  "
$ws.Range("D32").Value = "'7.12"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("D33").Value = "'622.25"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "'12.04"
$ws.Range("E34").Value = "  -5.56%  "
$ws.Range("D35").Value = "'63.20"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("E36").Value = "  -9.14%  "
$ws.Range("D37").Value = "'40.98"
$ws.Range("E37").Value = "  -10.86%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'0.392"
$ws.Range("E39").Value = "  -5.56%  "
$ws.Range("D40").Value = "0.0₃0745"
$ws.Range("E40").Value = "  -9.43%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -7.65%  "
$ws.Range("D43").Value = "2.946.67"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("D44").Value = "'2.77"
$ws.Range("E44").Value = "  -9.21%  "
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = "  -6.28%  "
$ws.Range("D46").Value = "'0.0400"
$ws.Range("E46").Value = "  -10.22%  "
$ws.Range("D47").Value = "'3.13"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  -8.00%  "
$ws.Range("D49").Value = "'138.14"
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("D50").Value = "'8.34"
$ws.Range("E50").Value = "  -8.36%  "
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'2.48"
$ws.Range("E51").Value = "  -9.35%  "
